# IWP30 test data refresh (Katalon RAD data sheet):
#  - Rows 5, 10, 11 and 12 were re-executed: their Date (col B) gets a
#    fresh run timestamp, Execute (col C) stays "Y".
#  - All other data rows (2,3,4,6,7,8,9,13,14) are left out of this run:
#    their Date stays as-is, but the Execute cell is cleared entirely so
#    Katalon skips them next time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New run timestamps for the re-executed rows only.
$ws.Range("B5").Value  = "Tue Jun 13 13:50:51 EDT 2023"
$ws.Range("B10").Value = "Tue Jun 13 13:51:40 EDT 2023"
$ws.Range("B11").Value = "Tue Jun 13 13:52:23 EDT 2023"
$ws.Range("B12").Value = "Tue Jun 13 13:53:05 EDT 2023"

# Clear the Execute flag (column C) for the rows skipped this pass -- the
# cell is removed outright, not just blanked.
$ws.Range("C2").Clear()
$ws.Range("C3").Clear()
$ws.Range("C4").Clear()
$ws.Range("C6").Clear()
$ws.Range("C7").Clear()
$ws.Range("C8").Clear()
$ws.Range("C9").Clear()
$ws.Range("C13").Clear()
$ws.Range("C14").Clear()

# Rows 5, 10, 11, 12 keep their Execute="Y" flag untouched.

# Match the author's final selection/scroll position.
$ws.Range("C13:C14").Select() | Out-Null
